# Week9 ObjectOrientedProgramming.pptx - minor lecture update.
#
# The code sample on the "main()" slide used pointer-member-access syntax
# (carX->method(...)) even though car1/car2 are plain (non-pointer)
# Automobile objects. Fix it to use dot-member-access (carX.method(...)).
#
# Slide 10 ("Title 1" placeholder) holds the whole code listing as one
# paragraph made of many same-language/same-format <a:r> runs separated by
# <a:br/>. We locate each "carX->..." occurrence inside the flattened
# TextRange.Text and fix it in place via TextRange.Characters(start,length),
# which preserves each run's existing character formatting (color/font/size).

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(10)
$shp = $s.Shapes.Item(1)
$tr  = $shp.TextFrame.TextRange

# Plain "->" -> "." fixes; each keeps the whole original run (including its
# leading 4-space indent) together as a single run, exactly like the source
# run it was already part of.
$simplePatterns = @(
    "    car1->setFuelEfficiency(8.2);",
    "    car1->drive(200.0);",
    "    car2->setFuelEfficiency(7.8);",
    "    car2->drive(200.0);",
    "    car1->setFuelEfficiency(6.2);",
    "    car1->drive(300.0);",
    "    car2->setFuelEfficiency(5.8);",
    "    car2->drive(300.0);",
    "    car1->displayReport();",
    "    car2->displayReport();"
)

# The two "addFuel" lines got edited in two separate passes by the author,
# which leaves the statement split across two runs: "...carX.addFuel" and
# the trailing "(50.0);" (both keep identical character formatting).
$splitPatterns = @(
    @{ old = "    car1->addFuel(50.0);"; part1 = "    car1.addFuel"; part2 = "(50.0);" },
    @{ old = "    car2->addFuel(50.0);"; part1 = "    car2.addFuel"; part2 = "(50.0);" }
)

# Resolve every edit's current offset up front, then apply them from the
# right-most offset to the left-most one, so earlier offsets are never
# invalidated by a preceding edit shrinking the text.
$edits = New-Object System.Collections.ArrayList

foreach ($pat in $simplePatterns) {
    $idx = $tr.Text.IndexOf($pat)
    if ($idx -ge 0) {
        [void]$edits.Add(@{ start = $idx; kind = "simple"; old = $pat })
    }
}
foreach ($sp in $splitPatterns) {
    $idx = $tr.Text.IndexOf($sp.old)
    if ($idx -ge 0) {
        [void]$edits.Add(@{ start = $idx; kind = "split"; old = $sp.old; part1 = $sp.part1; part2 = $sp.part2 })
    }
}

$orderedEdits = $edits | Sort-Object -Property start -Descending

foreach ($e in $orderedEdits) {
    if ($e.kind -eq "simple") {
        $newText = $e.old.Replace("->", ".")
        $rng = $tr.Characters($e.start + 1, $e.old.Length)
        $rng.Text = $newText
    } else {
        $rng1 = $tr.Characters($e.start + 1, $e.old.Length - $e.part2.Length)
        $rng1.Text = $e.part1

        # The text shrank by one character ("->" -> "."), so re-find the
        # trailing "(50.0);" piece rather than trusting the old offset math.
        $idx2 = $tr.Text.IndexOf($e.part2, $e.start)
        $rng2 = $tr.Characters($idx2 + 1, $e.part2.Length)
        $rng2.Text = $e.part2
    }
}
